# This workbook documents a FHIR StructureDefinition ("snapshot-age-in-years")
# as two sheets:
#   - "Metadata": a Property/Value key-value table
#   - "Elements": the element table (Path, Short, Definition, ...)
#
# The commit bumps the IG version/date, fills in the Publisher/Jurisdiction
# metadata (replacing a stray duplicated "Contact" row), and gives the root
# Extension element a real Short/Definition instead of the generic
# "Extension" / "An Extension" placeholder text.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$metaWs = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$metaWs.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$metaWs.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now "Alvearie Team"
$metaWs.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a second "Contact" / "No display for ContactDetail"
# row (a duplicate of row 11). It becomes the real Jurisdiction row.
$metaWs.Range("A10").Value = "Jurisdiction"
$metaWs.Range("B10").Value = "United States of America"

# Row 11 (the original duplicate "Contact" row) is removed outright,
# shifting every following row (Description, Purpose, ... Context) up by one.
$metaWs.Rows.Item(11).Delete()

# ---- Elements sheet ---------------------------------------------------
$elemWs = $wb.Worksheets.Item("Elements")

# Root "Extension" element (row 2): give it a real Short/Definition
# instead of the boilerplate "Extension" / "An Extension" text.
$elemWs.Range("K2").Value = "Age In Years"
$elemWs.Range("L2").Value = "Age in years at the time of the event"
